# Work Profile and new tenant support
# Appends new sprint-run history rows to the "AMSIN" and "AMS" sheets, and
# fixes up the formatting on AMS!A24:G24 (which was missing its normal
# cell style before this edit).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: write one data row (Run Date, Run Time, Sprint Name, Total,
# Pass, Fail, Time Taken) into columns A:G of $row on worksheet $ws.
#
# - Column A holds a text date like "2023-03-09". Excel's COM layer
#   auto-converts such strings into date serials on assignment, so we
#   prefix with an apostrophe to force it to stay literal text (matches
#   the workbook's existing inlineStr cells for this column).
# - Column B holds the numeric run-time serial, formatted with the
#   sheet's custom timestamp format so it still displays/round-trips as
#   a date-time like the existing rows.
# - Columns C..G are plain text/number values.
# - $normalize, when $true, forces the cell style back to the sheet's
#   plain "Normal" style (used for the AMS sheet, which has no
#   column-level default style so fresh cells come back completely
#   unstyled otherwise).
# ---------------------------------------------------------------------
function Set-HistoryRow {
    param(
        [object]$ws,
        [int]$row,
        [string]$runDate,
        [double]$runTime,
        [string]$sprintName,
        [double]$total,
        [double]$pass,
        [double]$fail,
        [double]$taken,
        [bool]$normalize
    )

    $cA = $ws.Cells.Item($row, 1)
    $cA.Value2 = "'" + $runDate
    if ($normalize) { $cA.Style = "Normal" }

    $cB = $ws.Cells.Item($row, 2)
    $cB.Value2 = $runTime
    $cB.NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $cC = $ws.Cells.Item($row, 3)
    $cC.Value2 = $sprintName
    if ($normalize) { $cC.Style = "Normal" }

    $cD = $ws.Cells.Item($row, 4)
    $cD.Value2 = $total
    if ($normalize) { $cD.Style = "Normal" }

    $cE = $ws.Cells.Item($row, 5)
    $cE.Value2 = $pass
    if ($normalize) { $cE.Style = "Normal" }

    $cF = $ws.Cells.Item($row, 6)
    $cF.Value2 = $fail
    if ($normalize) { $cF.Style = "Normal" }

    $cG = $ws.Cells.Item($row, 7)
    $cG.Value2 = $taken
    if ($normalize) { $cG.Style = "Normal" }
}

# ---------------------------------------------------------------------
# AMSIN sheet: rows 37-43 appended (dimension A1:G36 -> A1:G43).
# Columns default to style "5" on this sheet (set at the <cols> level),
# so new cells inherit it automatically and need no extra styling call
# - except for the very last row (43), whose cells come through with no
#   explicit style at all, matching the source diff.
# ---------------------------------------------------------------------
$wsAmsin = $wb.Worksheets.Item("AMSIN")

Set-HistoryRow $wsAmsin 37 "2023-03-09" 44994.5634352199    "174educfstcyle" 60 60 0 1    $false
Set-HistoryRow $wsAmsin 38 "2023-03-13" 44998.44833234954   "174eduflow"     60 60 0 1.09  $false
Set-HistoryRow $wsAmsin 39 "2023-03-30" 45015.70149211805   "175edusc"       60 59 1 1.15  $false
Set-HistoryRow $wsAmsin 40 "2023-03-31" 45016.49520108796   "175fnledu"      60 59 1 0.87  $false
Set-HistoryRow $wsAmsin 41 "2023-03-31" 45016.71950912037   "sadd"           59 59 0 1.43  $false
Set-HistoryRow $wsAmsin 42 "2023-04-06" 45022.64281496528   "176neweduc"     59 57 2 2.56  $false

# Row 43 keeps no explicit style on A/C/D/E/F/G (only B carries the
# date-time format) - force everything back to "Normal" to strip the
# column-inherited style.
Set-HistoryRow $wsAmsin 43 "2023-04-07" 45023.66776613814   "176fstedu"      59 59 0 0.84  $true
$wsAmsin.Cells.Item(43, 1).Style = "Normal"
$wsAmsin.Cells.Item(43, 3).Style = "Normal"

# ---------------------------------------------------------------------
# AMS sheet: rows 25-29 appended (dimension A1:G24 -> A1:G29), plus a
# formatting/value touch-up on the existing row 24, which previously had
# no explicit cell style and a very slightly different run-time value.
# This sheet has no column-level default style, so every new/changed
# cell is explicitly normalized.
# ---------------------------------------------------------------------
$wsAms = $wb.Worksheets.Item("AMS")

# --- fix up existing row 24 ---
$wsAms.Cells.Item(24, 1).Style = "Normal"
$wsAms.Cells.Item(24, 2).Value2 = 44977.82648603009
$wsAms.Cells.Item(24, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsAms.Cells.Item(24, 3).Style = "Normal"
$wsAms.Cells.Item(24, 4).Style = "Normal"
$wsAms.Cells.Item(24, 5).Style = "Normal"
$wsAms.Cells.Item(24, 6).Style = "Normal"
$wsAms.Cells.Item(24, 7).Style = "Normal"

# --- new rows ---
Set-HistoryRow $wsAms 25 "2023-03-13" 44998.55685114583 "174betaeduc" 60 60 0 0.96               $true
Set-HistoryRow $wsAms 26 "2023-03-13" 44998.85708618056 "174liveeedu" 60 60 0 0.86               $true
Set-HistoryRow $wsAms 27 "2023-03-31" 45016.56578296296 "175bted"     60 60 0 2.2                $true
Set-HistoryRow $wsAms 28 "2023-03-31" 45016.72568778935 "175edy"      59 59 0 0.8100000000000001 $true
Set-HistoryRow $wsAms 29 "2023-03-31" 45016.81429121528 "175liveedu"  59 59 0 0.76               $true
